$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet3")

# "Human" has been renamed "Homo sapiens (Human)" everywhere it appears on this sheet.
$ws.Range("F3").Value = "Homo sapiens (Human)"
$ws.Range("F5").Value = "Homo sapiens (Human)"

$ws.Activate()
[void]$ws.Range("F5").Select()
